$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 3: "True" text -> boolean TRUE in G3
$ws.Range("G3").Value = $true

# Row 4: "4" text -> numeric 4 in G4
$ws.Range("G4").Value = 4

# Row 6: "True" text -> boolean TRUE in G6
$ws.Range("G6").Value = $true

# Update the active selection shown in the sheet view
$ws.Range("D18").Select()
